$d = $word.ActiveDocument
$vtab = [char]11

# ---------------------------------------------------------------------------
# Paragraph 1 (currently the first picture): drop the picture, fill the
# paragraph with the text that used to sit in paragraph 2 - but with its
# final fragment corrected from "| -" to "| .".
# ---------------------------------------------------------------------------
$d.InlineShapes.Item(1).Delete()
$p1 = $d.Paragraphs.Item(1)
$text1 = "{a} Animal F and G have a prey-predator relationship. Based on the graph" + $vtab + "above} which animal is a prey and which is a predator? oe [1]" + $vtab + "| ."
$p1.Range.InsertBefore($text1)

# ---------------------------------------------------------------------------
# Paragraph 2 (was "{a} Animal F and G ..."): becomes "Preys |"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(
    "{a} Animal F and G have a prey-predator relationship. Based on the graph" + $vtab + "above} which animal is a prey and which is a predator? oe [1]" + $vtab + "| -",
    $true, $false, $false, $false, $false, $true, 1, $false, "Preys |", 2) | Out-Null

# ---------------------------------------------------------------------------
# Paragraph 3 (was "Prey: |"): becomes the "(b} Based en the graph ..." text.
# ---------------------------------------------------------------------------
$text3 = "(b} Based en the graph given above, there was @ period with very little" + $vtab + "rainfall! So"
$d.Content.Find.Execute("Prey: |", $true, $false, $false, $false, $false, $true, 1, $false, $text3, 2) | Out-Null

# ---------------------------------------------------------------------------
# Paragraph 4 (was "(b) Based on the graph ..."): removed outright so its
# successor ("This most likely happened ...") moves up to take its place.
# ---------------------------------------------------------------------------
$target4 = "(b) Based on the graph given above, there was a period with very little" + $vtab + "rainfall So"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq ($target4 + [char]13)) {
        $para.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------------
# Paragraph 5 (was "This most likely happened from Z fo! . 7]"):
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("This most likely happened from Z fo! . 7]", $true, $false, $false, $false, $false, $true, 1, $false, ". This most likely happened from Z fo | . {7", 2) | Out-Null

# ---------------------------------------------------------------------------
# Paragraph 6 ("(c) Other than ...") keeps its text, but the trailing empty
# paragraph and the second picture's paragraph must disappear. Word will
# never let the very last paragraph mark in the document be deleted, so
# instead of deleting paragraph 6 forward, we delete paragraph 6 and the
# following empty paragraph (both are *not* the final paragraph, so the
# deletions are legal), leaving the picture's paragraph as the new - and
# only remaining - final paragraph. We then strip its picture and restock
# it with paragraph 6's original text, which lets the forced-to-exist final
# paragraph mark serve as the document's true last paragraph.
# ---------------------------------------------------------------------------
$text6 = "(c) Other than lack of water for the predator, Give another possible reason. |" + $vtab + "why with very little. rainfall would cause a decrease in the population of" + $vtab + "the predator. = {1]"

$idx6 = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -eq ($text6 + [char]13)) {
        $idx6 = $i
        break
    }
}

if ($idx6 -gt 0) {
    # Delete paragraph 6 itself (text + mark) - legal, it is not the last paragraph.
    $d.Paragraphs.Item($idx6).Range.Delete()
    # Delete every following paragraph that is now empty and still not the last one.
    while ($idx6 -le $d.Paragraphs.Count -and $d.Paragraphs.Count -gt $idx6) {
        $p = $d.Paragraphs.Item($idx6)
        if ($p.Range.InlineShapes.Count -eq 0 -and $p.Range.Text -eq [char]13) {
            $p.Range.Delete()
        } else {
            break
        }
    }
    # Whatever paragraph remains at $idx6 is now the forced final paragraph of the
    # document (its picture, if any, is removed) - put paragraph 6's text back into it.
    $finalPara = $d.Paragraphs.Item($idx6)
    while ($finalPara.Range.InlineShapes.Count -gt 0) {
        $finalPara.Range.InlineShapes.Item(1).Delete()
    }
    $finalPara.Range.InsertBefore($text6)
}
